$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.1453803460577
$ws.Range("C2").Value = 11.09381386923738
$ws.Range("D2").Value = 4.991507198785488
$ws.Range("E2").Value = 12.33569114822264
$ws.Range("F2").Value = 24.49992230631182
$ws.Range("I2").Value = 22.19670225903991
$ws.Range("L2").Value = 10.03587379587387
$ws.Range("M2").Value = 14.53851040075913
$ws.Range("O2").Value = 21.94641149117582
$ws.Range("B3").Value = 14.54781237916831
$ws.Range("C3").Value = 10.81483307787347
$ws.Range("D3").Value = 4.94668022060366
$ws.Range("E3").Value = 12.38353142493918
$ws.Range("F3").Value = 24.55419685818636
$ws.Range("I3").Value = 22.34220659126362
$ws.Range("L3").Value = 10.04495477195745
$ws.Range("M3").Value = 14.40417993689605
$ws.Range("O3").Value = 22.04265431265349
$ws.Range("B4").Value = 14.16874570717541
$ws.Range("C4").Value = 10.63888572589539
$ws.Range("D4").Value = 4.91874607288396
$ws.Range("E4").Value = 12.41445368153827
$ws.Range("F4").Value = 24.59586488378161
$ws.Range("I4").Value = 22.43698940105832
$ws.Range("L4").Value = 10.05196209405618
$ws.Range("M4").Value = 14.32267194227352
$ws.Range("O4").Value = 22.10802811107588
$ws.Range("B5").Value = 14.0114259447601
$ws.Range("C5").Value = 10.56608979466771
$ws.Range("D5").Value = 4.907265396023574
$ws.Range("E5").Value = 12.42744520871785
$ws.Range("F5").Value = 24.61493434888962
$ws.Range("I5").Value = 22.47698246070599
$ws.Range("L5").Value = 10.05517801652476
$ws.Range("M5").Value = 14.28972932047398
$ws.Range("O5").Value = 22.13624164720366
$ws.Range("B6").Value = 13.98513779574668
$ws.Range("C6").Value = 10.553938145621
$ws.Range("D6").Value = 4.90535333472902
$ws.Range("E6").Value = 12.42962605882853
$ws.Range("F6").Value = 24.61822674386162
$ws.Range("I6").Value = 22.48370591301035
$ws.Range("L6").Value = 10.0557337946026
$ws.Range("M6").Value = 14.28427650101864
$ws.Range("O6").Value = 22.14102132014259
$ws.Range("B7").Value = 14.16663526937799
$ws.Range("C7").Value = 10.63790831086821
$ws.Range("D7").Value = 4.918591626598002
$ws.Range("E7").Value = 12.41462730732299
$ws.Range("F7").Value = 24.59611361342521
$ws.Range("I7").Value = 22.43752322159789
$ws.Range("L7").Value = 10.052004005457
$ws.Range("M7").Value = 14.3222265260409
$ws.Range("O7").Value = 22.10840224705852
$ws.Range("B8").Value = 14.94198064445215
$ws.Range("C8").Value = 10.99862959271919
$ws.Range("D8").Value = 4.976139494802282
$ws.Range("E8").Value = 12.35186586394249
$ws.Range("F8").Value = 24.51689968799357
$ws.Range("I8").Value = 22.24574236534089
$ws.Range("L8").Value = 10.03870805523311
$ws.Range("M8").Value = 14.4920087230655
$ws.Range("O8").Value = 21.97828893800851
$ws.Range("B9").Value = 16.35795579589243
$ws.Range("C9").Value = 11.66610024122555
$ws.Range("D9").Value = 5.08548238843728
$ws.Range("E9").Value = 12.24102123337154
$ws.Range("F9").Value = 24.42809738405583
$ws.Range("I9").Value = 21.91285814054458
$ws.Range("L9").Value = 10.0239740412699
$ws.Range("M9").Value = 14.83140026703189
$ws.Range("O9").Value = 21.77323334470514
$ws.Range("B10").Value = 17.32542332544563
$ws.Range("C10").Value = 12.12851602970492
$ws.Range("D10").Value = 5.163344934550923
$ws.Range("E10").Value = 12.1669646859639
$ws.Range("F10").Value = 24.40380503533924
$ws.Range("I10").Value = 21.69464176256769
$ws.Range("L10").Value = 10.02003269903116
$ws.Range("M10").Value = 15.08298202964175
$ws.Range("O10").Value = 21.65347172931441
$ws.Range("B11").Value = 17.74819780227867
$ws.Range("C11").Value = 12.3321343486555
$ws.Range("D11").Value = 5.198161007245979
$ws.Range("E11").Value = 12.13486120757229
$ws.Range("F11").Value = 24.40170436372912
$ws.Range("I11").Value = 21.60109831678488
$ws.Range("L11").Value = 10.01972676734745
$ws.Range("M11").Value = 15.19755475806755
$ws.Range("O11").Value = 21.605769212011
$ws.Range("B12").Value = 17.90569467384247
$ws.Range("C12").Value = 12.408220856067
$ws.Range("D12").Value = 5.211252402036338
$ws.Range("E12").Value = 12.12293120608059
$ws.Range("F12").Value = 24.40219862726384
$ws.Range("I12").Value = 21.56649994658942
$ws.Range("L12").Value = 10.01982394753599
$ws.Range("M12").Value = 15.24092842233447
$ws.Range("O12").Value = 21.58868577796178
$ws.Range("B13").Value = 17.87189201278609
$ws.Range("C13").Value = 12.39188035618757
$ws.Range("D13").Value = 5.20843715027572
$ws.Range("E13").Value = 12.12549047124349
$ws.Range("F13").Value = 24.40203478763351
$ws.Range("I13").Value = 21.57391463169346
$ws.Range("L13").Value = 10.01979355709598
$ws.Range("M13").Value = 15.23158817200982
$ws.Range("O13").Value = 21.59232130345647
$ws.Range("B14").Value = 17.76120780072247
$ws.Range("C14").Value = 12.33841473345342
$ws.Range("D14").Value = 5.199239934953813
$ws.Range("E14").Value = 12.1338751780338
$ws.Range("F14").Value = 24.40171916775638
$ws.Range("I14").Value = 21.59823535555532
$ws.Range("L14").Value = 10.01973049785775
$ws.Range("M14").Value = 15.20112356201832
$ws.Range("O14").Value = 21.60434406471
$ws.Range("B15").Value = 17.69306920776175
$ws.Range("C15").Value = 12.3055313274298
$ws.Range("D15").Value = 5.193594131999339
$ws.Range("E15").Value = 12.13904056929291
$ws.Range("F15").Value = 24.40169386266032
$ws.Range("I15").Value = 21.61323991821912
$ws.Range("L15").Value = 10.0197195898902
$ws.Range("M15").Value = 15.18246056745505
$ws.Range("O15").Value = 21.61183620975244
$ws.Range("B16").Value = 17.2974358157285
$ws.Range("C16").Value = 12.11506909982274
$ws.Range("D16").Value = 5.16105696607658
$ws.Range("E16").Value = 12.16909453130305
$ws.Range("F16").Value = 24.40412268095938
$ws.Range("I16").Value = 21.70087034708135
$ws.Range("L16").Value = 10.02008254825901
$ws.Range("M16").Value = 15.07549434672858
$ws.Range("O16").Value = 21.65672609030289
$ws.Range("B17").Value = 17.05020785426925
$ws.Range("C17").Value = 11.99646373021664
$ws.Range("D17").Value = 5.140937574740219
$ws.Range("E17").Value = 12.18793691627085
$ws.Range("F17").Value = 24.407907367413
$ws.Range("I17").Value = 21.75609582839667
$ws.Range("L17").Value = 10.02068552634801
$ws.Range("M17").Value = 15.00988453547242
$ws.Range("O17").Value = 21.68600473278964
$ws.Range("B18").Value = 16.9063828699717
$ws.Range("C18").Value = 11.92761366492989
$ws.Range("D18").Value = 5.129308924387944
$ws.Range("E18").Value = 12.1989238288203
$ws.Range("F18").Value = 24.41092655797131
$ws.Range("I18").Value = 21.78839877167468
$ws.Range("L18").Value = 10.02117236036835
$ws.Range("M18").Value = 14.97216127987251
$ws.Range("O18").Value = 21.70348264058424
$ws.Range("B19").Value = 16.85741043344856
$ws.Range("C19").Value = 11.90419539091005
$ws.Range("D19").Value = 5.125362137911709
$ws.Range("E19").Value = 12.20266947862097
$ws.Range("F19").Value = 24.41209337161257
$ws.Range("I19").Value = 21.79942847740902
$ws.Range("L19").Value = 10.02136126506296
$ws.Range("M19").Value = 14.95939214259243
$ws.Range("O19").Value = 21.70950968738884
$ws.Range("B20").Value = 17.07669479760098
$ws.Range("C20").Value = 12.0091552022615
$ws.Range("D20").Value = 5.143085209449005
$ws.Range("E20").Value = 12.18591567095613
$ws.Range("F20").Value = 24.4074172781124
$ws.Range("I20").Value = 21.75016121801048
$ws.Range("L20").Value = 10.02060685260022
$ws.Range("M20").Value = 15.01686761057632
$ws.Range("O20").Value = 21.68282193821559
$ws.Range("B21").Value = 17.79378973575117
$ws.Range("C21").Value = 12.35414692718705
$ws.Range("D21").Value = 5.20194394171131
$ws.Range("E21").Value = 12.13140623712737
$ws.Range("F21").Value = 24.40177685447382
$ws.Range("I21").Value = 21.5910693811938
$ws.Range("L21").Value = 10.01974324509751
$ws.Range("M21").Value = 15.21007234049054
$ws.Range("O21").Value = 21.60078603295712
$ws.Range("B22").Value = 18.24726689532443
$ws.Range("C22").Value = 12.57365729951382
$ws.Range("D22").Value = 5.23986860447826
$ws.Range("E22").Value = 12.09710316684529
$ws.Range("F22").Value = 24.40560844671403
$ws.Range("I22").Value = 21.49190022235826
$ws.Range("L22").Value = 10.02042009936728
$ws.Range("M22").Value = 15.33625568715754
$ws.Range("O22").Value = 21.55288833282459
$ws.Range("B23").Value = 18.00665839134718
$ws.Range("C23").Value = 12.45706149036353
$ws.Range("D23").Value = 5.219679071771859
$ws.Range("E23").Value = 12.1152907473347
$ws.Range("F23").Value = 24.40287500108353
$ws.Range("I23").Value = 21.54438838692634
$ws.Range("L23").Value = 10.0199455586241
$ws.Range("M23").Value = 15.26892730825576
$ws.Range("O23").Value = 21.57792720897513
$ws.Range("B24").Value = 17.06472530731554
$ws.Range("C24").Value = 12.00341944669179
$ws.Range("D24").Value = 5.142114455664979
$ws.Range("E24").Value = 12.18682899564037
$ws.Range("F24").Value = 24.40763622079594
$ws.Range("I24").Value = 21.75284253262479
$ws.Range("L24").Value = 10.02064198431525
$ws.Range("M24").Value = 15.01371057024489
$ws.Range("O24").Value = 21.68425887022205
$ws.Range("B25").Value = 15.98710838541664
$ws.Range("C25").Value = 11.49019986289859
$ws.Range("D25").Value = 5.05631250487631
$ws.Range("E25").Value = 12.26970607836615
$ws.Range("F25").Value = 24.44495288600158
$ws.Range("I25").Value = 21.99828599568501
$ws.Range("L25").Value = 10.02674872057311
$ws.Range("M25").Value = 14.73907054651487
$ws.Range("O25").Value = 21.82330608012984
